$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (ALC)
$ws.Cells.Item(17, 8).Value = 2395224.2  # H17: 3913604.5 -> 2395224.2
$ws.Cells.Item(17, 10).Value = 2431508  # J17: 4011432.2 -> 2431508
$ws.Cells.Item(17, 12).Value = 7294524  # L17: 12034296.6 -> 7294524
$ws.Cells.Item(17, 14).Value = -7294860  # N17: -12034632.6 -> -7294860

# Row 19 (ALC)
$ws.Cells.Item(19, 8).Value = 473  # H19: 647.5 -> 473
$ws.Cells.Item(19, 9).Value = 277.14285  # I19: 500.33334 -> 277.14285
$ws.Cells.Item(19, 10).Value = 930  # J19: 735.8 -> 930
$ws.Cells.Item(19, 11).Value = 277.14285  # K19: 500.33334 -> 277.14285
$ws.Cells.Item(19, 12).Value = 930  # L19: 735.8 -> 930
$ws.Cells.Item(19, 13).Value = -102.14285  # M19: -325.33334 -> -102.14285
$ws.Cells.Item(19, 14).Value = -1280  # N19: -1085.8 -> -1280

# Row 116 (ALC)
$ws.Cells.Item(116, 8).Value = 2274.0435  # H116: 2424.25 -> 2274.0435
$ws.Cells.Item(116, 9).Value = 1592.5385  # I116: 1645 -> 1592.5385
$ws.Cells.Item(116, 10).Value = 3160  # J116: 3061.818 -> 3160
$ws.Cells.Item(116, 11).Value = 1592.5385  # K116: 1645 -> 1592.5385
$ws.Cells.Item(116, 12).Value = 3160  # L116: 3061.818 -> 3160
$ws.Cells.Item(116, 13).Value = 1849.4615  # M116: 1797 -> 1849.4615
$ws.Cells.Item(116, 14).Value = -10044  # N116: -9945.817999999999 -> -10044

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Cells.Item(2, 8).Value = 1981.3334  # H2: 1714.4 -> 1981.3334
$ws.Cells.Item(2, 9).Value = 1847  # I2: 1724.8334 -> 1847
$ws.Cells.Item(2, 10).Value = 2250  # J2: 1698.75 -> 2250
$ws.Cells.Item(2, 11).Value = 1847  # K2: 1724.8334 -> 1847
$ws.Cells.Item(2, 12).Value = 2250  # L2: 1698.75 -> 2250
$ws.Cells.Item(2, 13).Value = -1734  # M2: -1611.8334 -> -1734
$ws.Cells.Item(2, 14).Value = -2476  # N2: -1924.75 -> -2476

# Row 74 (ARM)
$ws.Cells.Item(74, 8).Value = 1708.409  # H74: 1196 -> 1708.409
$ws.Cells.Item(74, 9).Value = 1455.4706  # I74: 967.26666 -> 1455.4706
$ws.Cells.Item(74, 11).Value = 1455.4706  # K74: 967.26666 -> 1455.4706
$ws.Cells.Item(74, 13).Value = -581.4706000000001  # M74: -93.26666 -> -581.4706000000001

# Row 77 (ARM)
$ws.Cells.Item(77, 8).Value = 1708.409  # H77: 1196 -> 1708.409
$ws.Cells.Item(77, 9).Value = 1455.4706  # I77: 967.26666 -> 1455.4706
$ws.Cells.Item(77, 11).Value = 7277.353000000001  # K77: 4836.3333 -> 7277.353000000001
$ws.Cells.Item(77, 13).Value = -2909.353000000001  # M77: -468.3333000000002 -> -2909.353000000001

# Row 105 (ARM)
$ws.Cells.Item(105, 8).Value = 28000  # H105: 0 -> 28000
$ws.Cells.Item(105, 10).Value = 28000  # J105: 0 -> 28000
$ws.Cells.Item(105, 12).Value = 28000  # L105: 0 -> 28000
$ws.Cells.Item(105, 14).Value = -34988  # N105: None -> -34988

# Row 116 (ARM)
$ws.Cells.Item(116, 8).Value = 1981.3334  # H116: 1714.4 -> 1981.3334
$ws.Cells.Item(116, 9).Value = 1847  # I116: 1724.8334 -> 1847
$ws.Cells.Item(116, 10).Value = 2250  # J116: 1698.75 -> 2250
$ws.Cells.Item(116, 11).Value = 1847  # K116: 1724.8334 -> 1847
$ws.Cells.Item(116, 12).Value = 2250  # L116: 1698.75 -> 2250
$ws.Cells.Item(116, 13).Value = 447  # M116: 569.1666 -> 447
$ws.Cells.Item(116, 14).Value = -6838  # N116: -6286.75 -> -6838

# Row 132 (ARM)
$ws.Cells.Item(132, 8).Value = 3675.1482  # H132: 3422.3 -> 3675.1482
$ws.Cells.Item(132, 9).Value = 2202.9092  # I132: 1976.5714 -> 2202.9092
$ws.Cells.Item(132, 11).Value = 6608.7276  # K132: 5929.7142 -> 6608.7276
$ws.Cells.Item(132, 13).Value = -4078.7276  # M132: -3399.7142 -> -4078.7276

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Cells.Item(3, 8).Value = 1981.3334  # H3: 1714.4 -> 1981.3334
$ws.Cells.Item(3, 9).Value = 1847  # I3: 1724.8334 -> 1847
$ws.Cells.Item(3, 10).Value = 2250  # J3: 1698.75 -> 2250
$ws.Cells.Item(3, 11).Value = 1847  # K3: 1724.8334 -> 1847
$ws.Cells.Item(3, 12).Value = 2250  # L3: 1698.75 -> 2250
$ws.Cells.Item(3, 13).Value = -1733  # M3: -1610.8334 -> -1733
$ws.Cells.Item(3, 14).Value = -2478  # N3: -1926.75 -> -2478

# Row 107 (BSM)
$ws.Cells.Item(107, 8).Value = 1074.8422  # H107: 1106.7778 -> 1074.8422
$ws.Cells.Item(107, 9).Value = 887.0909  # I107: 925.8 -> 887.0909
$ws.Cells.Item(107, 11).Value = 887.0909  # K107: 925.8 -> 887.0909
$ws.Cells.Item(107, 13).Value = 1032.9091  # M107: 994.2 -> 1032.9091

# Row 134 (BSM)
$ws.Cells.Item(134, 8).Value = 7924.35  # H134: 7914.35 -> 7924.35
$ws.Cells.Item(134, 9).Value = 12851.1  # I134: 12831.1 -> 12851.1
$ws.Cells.Item(134, 11).Value = 38553.3  # K134: 38493.3 -> 38553.3
$ws.Cells.Item(134, 13).Value = -36018.3  # M134: -35958.3 -> -36018.3

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (CRP)
$ws.Cells.Item(16, 8).Value = 2042.2142  # H16: 2273.818 -> 2042.2142
$ws.Cells.Item(16, 9).Value = 1775.3334  # I16: 2119.8 -> 1775.3334
$ws.Cells.Item(16, 10).Value = 2522.6  # J16: 2402.1667 -> 2522.6
$ws.Cells.Item(16, 11).Value = 1775.3334  # K16: 2119.8 -> 1775.3334
$ws.Cells.Item(16, 12).Value = 2522.6  # L16: 2402.1667 -> 2522.6
$ws.Cells.Item(16, 13).Value = -1488.3334  # M16: -1832.8 -> -1488.3334
$ws.Cells.Item(16, 14).Value = -3096.6  # N16: -2976.1667 -> -3096.6

# Row 58 (CRP)
$ws.Cells.Item(58, 8).Value = 1680.4324  # H58: 2180.5625 -> 1680.4324
$ws.Cells.Item(58, 9).Value = 1201.6842  # I58: 1621.2727 -> 1201.6842
$ws.Cells.Item(58, 10).Value = 2185.7778  # J58: 2473.524 -> 2185.7778
$ws.Cells.Item(58, 11).Value = 1201.6842  # K58: 1621.2727 -> 1201.6842
$ws.Cells.Item(58, 12).Value = 2185.7778  # L58: 2473.524 -> 2185.7778
$ws.Cells.Item(58, 13).Value = -998.6841999999999  # M58: -1418.2727 -> -998.6841999999999
$ws.Cells.Item(58, 14).Value = -2591.7778  # N58: -2879.524 -> -2591.7778

# Row 94 (CRP)
$ws.Cells.Item(94, 8).Value = 4633.9585  # H94: 3838.375 -> 4633.9585
$ws.Cells.Item(94, 9).Value = 4487.5557  # I94: 2939.7334 -> 4487.5557
$ws.Cells.Item(94, 10).Value = 4721.8  # J94: 4631.294 -> 4721.8
$ws.Cells.Item(94, 11).Value = 4487.5557  # K94: 2939.7334 -> 4487.5557
$ws.Cells.Item(94, 12).Value = 4721.8  # L94: 4631.294 -> 4721.8
$ws.Cells.Item(94, 13).Value = -4036.5557  # M94: -2488.7334 -> -4036.5557
$ws.Cells.Item(94, 14).Value = -5623.8  # N94: -5533.294 -> -5623.8

# Row 99 (CRP)
$ws.Cells.Item(99, 8).Value = 1051.25  # H99: 1064.4445 -> 1051.25
$ws.Cells.Item(99, 9).Value = 902.2  # I99: 963.3333 -> 902.2
$ws.Cells.Item(99, 10).Value = 1299.6666  # J99: 1266.6666 -> 1299.6666
$ws.Cells.Item(99, 11).Value = 902.2  # K99: 963.3333 -> 902.2
$ws.Cells.Item(99, 12).Value = 1299.6666  # L99: 1266.6666 -> 1299.6666
$ws.Cells.Item(99, 13).Value = 595.8  # M99: 534.6667 -> 595.8
$ws.Cells.Item(99, 14).Value = -4295.6666  # N99: -4262.6666 -> -4295.6666

# Row 113 (CRP)
$ws.Cells.Item(113, 8).Value = 2042.2142  # H113: 2273.818 -> 2042.2142
$ws.Cells.Item(113, 9).Value = 1775.3334  # I113: 2119.8 -> 1775.3334
$ws.Cells.Item(113, 10).Value = 2522.6  # J113: 2402.1667 -> 2522.6
$ws.Cells.Item(113, 11).Value = 1775.3334  # K113: 2119.8 -> 1775.3334
$ws.Cells.Item(113, 12).Value = 2522.6  # L113: 2402.1667 -> 2522.6
$ws.Cells.Item(113, 13).Value = 394.6666  # M113: 50.19999999999982 -> 394.6666
$ws.Cells.Item(113, 14).Value = -6862.6  # N113: -6742.1667 -> -6862.6

# Row 122 (CRP)
$ws.Cells.Item(122, 8).Value = 13898888  # H122: 3904 -> 13898888
$ws.Cells.Item(122, 9).Value = 27777776  # I122: 303.66666 -> 27777776
$ws.Cells.Item(122, 10).Value = 20000  # J122: 7504.3335 -> 20000
$ws.Cells.Item(122, 11).Value = 83333328  # K122: 910.9999799999999 -> 83333328
$ws.Cells.Item(122, 12).Value = 60000  # L122: 22513.0005 -> 60000
$ws.Cells.Item(122, 13).Value = -83330878  # M122: 1539.00002 -> -83330878
$ws.Cells.Item(122, 14).Value = -64900  # N122: -27413.0005 -> -64900

# Row 126 (CRP)
$ws.Cells.Item(126, 8).Value = 1051.25  # H126: 1064.4445 -> 1051.25
$ws.Cells.Item(126, 9).Value = 902.2  # I126: 963.3333 -> 902.2
$ws.Cells.Item(126, 10).Value = 1299.6666  # J126: 1266.6666 -> 1299.6666
$ws.Cells.Item(126, 11).Value = 2706.6  # K126: 2889.9999 -> 2706.6
$ws.Cells.Item(126, 12).Value = 3898.9998  # L126: 3799.9998 -> 3898.9998
$ws.Cells.Item(126, 13).Value = -236.6000000000004  # M126: -419.9998999999998 -> -236.6000000000004
$ws.Cells.Item(126, 14).Value = -8838.9998  # N126: -8739.9998 -> -8838.9998

# Row 132 (CRP)
$ws.Cells.Item(132, 8).Value = 1831.2759  # H132: 2268.8572 -> 1831.2759
$ws.Cells.Item(132, 9).Value = 1522.2142  # I132: 1994.4445 -> 1522.2142
$ws.Cells.Item(132, 10).Value = 2119.7334  # J132: 2474.6667 -> 2119.7334
$ws.Cells.Item(132, 11).Value = 4566.642599999999  # K132: 5983.333500000001 -> 4566.642599999999
$ws.Cells.Item(132, 12).Value = 6359.2002  # L132: 7424.000100000001 -> 6359.2002
$ws.Cells.Item(132, 13).Value = -2036.642599999999  # M132: -3453.333500000001 -> -2036.642599999999
$ws.Cells.Item(132, 14).Value = -11419.2002  # N132: -12484.0001 -> -11419.2002

# Row 136 (CRP)
$ws.Cells.Item(136, 8).Value = 1680.4324  # H136: 2180.5625 -> 1680.4324
$ws.Cells.Item(136, 9).Value = 1201.6842  # I136: 1621.2727 -> 1201.6842
$ws.Cells.Item(136, 10).Value = 2185.7778  # J136: 2473.524 -> 2185.7778
$ws.Cells.Item(136, 11).Value = 3605.0526  # K136: 4863.8181 -> 3605.0526
$ws.Cells.Item(136, 12).Value = 6557.3334  # L136: 7420.572 -> 6557.3334
$ws.Cells.Item(136, 13).Value = -1055.0526  # M136: -2313.8181 -> -1055.0526
$ws.Cells.Item(136, 14).Value = -11657.3334  # N136: -12520.572 -> -11657.3334

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (CUL)
$ws.Cells.Item(5, 8).Value = 218324.05  # H5: 218328.06 -> 218324.05
$ws.Cells.Item(5, 9).Value = 463.44446  # I5: 480.82352 -> 463.44446
$ws.Cells.Item(5, 10).Value = 358377.28  # J5: 346031.62 -> 358377.28
$ws.Cells.Item(5, 11).Value = 1390.33338  # K5: 1442.47056 -> 1390.33338
$ws.Cells.Item(5, 12).Value = 1075131.84  # L5: 1038094.86 -> 1075131.84
$ws.Cells.Item(5, 13).Value = -1278.33338  # M5: -1330.47056 -> -1278.33338
$ws.Cells.Item(5, 14).Value = -1075355.84  # N5: -1038318.86 -> -1075355.84

# Row 80 (CUL)
$ws.Cells.Item(80, 8).Value = 2033.3334  # H80: 2636.7646 -> 2033.3334
$ws.Cells.Item(80, 9).Value = 1550  # I80: 2100 -> 1550
$ws.Cells.Item(80, 10).Value = 3000  # J80: 2708.3333 -> 3000
$ws.Cells.Item(80, 11).Value = 4650  # K80: 6300 -> 4650
$ws.Cells.Item(80, 12).Value = 9000  # L80: 8124.999899999999 -> 9000
$ws.Cells.Item(80, 13).Value = -3714  # M80: -5364 -> -3714
$ws.Cells.Item(80, 14).Value = -10872  # N80: -9996.999899999999 -> -10872

# Row 83 (CUL)
$ws.Cells.Item(83, 8).Value = 2033.3334  # H83: 2636.7646 -> 2033.3334
$ws.Cells.Item(83, 9).Value = 1550  # I83: 2100 -> 1550
$ws.Cells.Item(83, 10).Value = 3000  # J83: 2708.3333 -> 3000
$ws.Cells.Item(83, 11).Value = 13950  # K83: 18900 -> 13950
$ws.Cells.Item(83, 12).Value = 27000  # L83: 24374.9997 -> 27000
$ws.Cells.Item(83, 13).Value = -9270  # M83: -14220 -> -9270
$ws.Cells.Item(83, 14).Value = -36360  # N83: -33734.9997 -> -36360

# Row 120 (CUL)
$ws.Cells.Item(120, 8).Value = 2307.5  # H120: 7116.6665 -> 2307.5
$ws.Cells.Item(120, 9).Value = 2307.5  # I120: 2233.3333 -> 2307.5
$ws.Cells.Item(120, 10).Value = 0  # J120: 12000 -> 0
$ws.Cells.Item(120, 11).Value = 6922.5  # K120: 6699.999899999999 -> 6922.5
$ws.Cells.Item(120, 12).Value = 0  # L120: 36000 -> 0
$ws.Cells.Item(120, 13).ClearContents()  # M120: -1861.999899999999 -> (removed)
$ws.Cells.Item(120, 14).Value = -2084.5  # N120: -45676 -> -2084.5

# Row 135 (CUL)
$ws.Cells.Item(135, 8).Value = 218324.05  # H135: 218328.06 -> 218324.05
$ws.Cells.Item(135, 9).Value = 463.44446  # I135: 480.82352 -> 463.44446
$ws.Cells.Item(135, 10).Value = 358377.28  # J135: 346031.62 -> 358377.28
$ws.Cells.Item(135, 11).Value = 4171.00014  # K135: 4327.41168 -> 4171.00014
$ws.Cells.Item(135, 12).Value = 3225395.52  # L135: 3114284.58 -> 3225395.52
$ws.Cells.Item(135, 13).Value = -1636.00014  # M135: -1792.41168 -> -1636.00014
$ws.Cells.Item(135, 14).Value = -3230465.52  # N135: -3119354.58 -> -3230465.52

$ws = $wb.Worksheets.Item("GSM")
# Row 126 (GSM)
$ws.Cells.Item(126, 8).Value = 4275.365  # H126: 4363.14 -> 4275.365
$ws.Cells.Item(126, 9).Value = 10085.333  # I126: 10071.833 -> 10085.333
$ws.Cells.Item(126, 10).Value = 2532.375  # J126: 2560.3948 -> 2532.375
$ws.Cells.Item(126, 11).Value = 30255.999  # K126: 30215.499 -> 30255.999
$ws.Cells.Item(126, 12).Value = 7597.125  # L126: 7681.1844 -> 7597.125
$ws.Cells.Item(126, 13).Value = -27785.999  # M126: -27745.499 -> -27785.999
$ws.Cells.Item(126, 14).Value = -12537.125  # N126: -12621.1844 -> -12537.125

# Row 132 (GSM)
$ws.Cells.Item(132, 8).Value = 2885.3  # H132: 2871.25 -> 2885.3
$ws.Cells.Item(132, 9).Value = 2998.1667  # I132: 2951.3333 -> 2998.1667
$ws.Cells.Item(132, 11).Value = 8994.500100000001  # K132: 8853.999899999999 -> 8994.500100000001
$ws.Cells.Item(132, 13).Value = -6464.500100000001  # M132: -6323.999899999999 -> -6464.500100000001

$ws = $wb.Worksheets.Item("LTW")
# Row 61 (LTW)
$ws.Cells.Item(61, 8).Value = 2056.6  # H61: 1452.5385 -> 2056.6
$ws.Cells.Item(61, 9).Value = 2122  # I61: 1424 -> 2122
$ws.Cells.Item(61, 11).Value = 2122  # K61: 1424 -> 2122
$ws.Cells.Item(61, 13).Value = -1920  # M61: -1222 -> -1920

# Row 93 (LTW)
$ws.Cells.Item(93, 8).Value = 852.2727  # H93: 959.6429000000001 -> 852.2727
$ws.Cells.Item(93, 9).Value = 627.1539  # I93: 703.8889 -> 627.1539
$ws.Cells.Item(93, 10).Value = 1177.4445  # J93: 1420 -> 1177.4445
$ws.Cells.Item(93, 11).Value = 627.1539  # K93: 703.8889 -> 627.1539
$ws.Cells.Item(93, 12).Value = 1177.4445  # L93: 1420 -> 1177.4445
$ws.Cells.Item(93, 13).Value = 620.8461  # M93: 544.1111 -> 620.8461
$ws.Cells.Item(93, 14).Value = -3673.4445  # N93: -3916 -> -3673.4445

# Row 113 (LTW)
$ws.Cells.Item(113, 8).Value = 2056.6  # H113: 1452.5385 -> 2056.6
$ws.Cells.Item(113, 9).Value = 2122  # I113: 1424 -> 2122
$ws.Cells.Item(113, 11).Value = 2122  # K113: 1424 -> 2122
$ws.Cells.Item(113, 13).Value = 48  # M113: 746 -> 48

# Row 122 (LTW)
$ws.Cells.Item(122, 8).Value = 3259947  # H122: 3018621.8 -> 3259947
$ws.Cells.Item(122, 9).Value = 5104484  # I122: 4466680.5 -> 5104484
$ws.Cells.Item(122, 11).Value = 15313452  # K122: 13400041.5 -> 15313452
$ws.Cells.Item(122, 13).Value = -15311002  # M122: -13397591.5 -> -15311002

# Row 136 (LTW)
$ws.Cells.Item(136, 8).Value = 7910.8047  # H136: 7910.878 -> 7910.8047
$ws.Cells.Item(136, 9).Value = 5469.467  # I136: 5469.567 -> 5469.467
$ws.Cells.Item(136, 11).Value = 16408.401  # K136: 16408.701 -> 16408.401
$ws.Cells.Item(136, 13).Value = -13858.401  # M136: -13858.701 -> -13858.401

$ws = $wb.Worksheets.Item("WVR")
# Row 136 (WVR)
$ws.Cells.Item(136, 8).Value = 1994.42  # H136: 2065.1042 -> 1994.42
$ws.Cells.Item(136, 9).Value = 2019.0333  # I136: 2141.9644 -> 2019.0333
$ws.Cells.Item(136, 11).Value = 6057.0999  # K136: 6425.8932 -> 6057.0999
$ws.Cells.Item(136, 13).Value = -3507.0999  # M136: -3875.8932 -> -3507.0999
Write-Host "All 34 row edits applied."
